# Apply cryptocurrency price/volume updates (Tue Sep 19 2023 GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of new Price values look numeric ("217.68", "0.0512", ...) but the
# column stores plain text (e.g. existing "1.00" cells), so force the Text number
# format first to stop Excel from auto-converting those assignments to numbers.

$ws.Range("D2").Value = "27.187.44"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.645.47"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.68"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.510"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0630"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.01"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "1.875.40"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "1.626.30"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.12"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.41"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "27.145.36"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.94"
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.44"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.50"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.18"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.59"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.42"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.119"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.75"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("D35").Value = "1.265.99"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.543"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.841"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +4.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.38"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "1.785.59"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.25"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.79"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +16.45%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.71"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0973"
$ws.Range("E51").Value = "  -1.16%  "
